$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: 19 Feb 2020 (W), 2200-0000 ---
$ws.Range("A36").Value = "19 Feb 2020 (W)"
$ws.Range("B36").Value = "2200-0000"
$ws.Range("C36").Value = "Harry, Deon, Thuc"
$ws.Range("D36").Value = "Supertask by working on latest homework, harassment training, and 264 lab at the same time"
$ws.Range("E36").Value = "Finished all, just need to review homework for this class with teammates tomorrow in lab"
$ws.Rows.Item(36).RowHeight = 36.5

# --- Row 37: continuation into 20 Feb 2020 (Th), 0000-0200 ---
$ws.Range("A37").Value = "20 Feb 2020 (Th)"
$ws.Range("B37").Value = "0000-0200"
$ws.Range("C37").Value = "Harry, Deon, Thuc"
$ws.Range("D37").Formula = "=D36"
$ws.Range("E37").Value = "Finished all, just need to review homework for this class with teammates today in lab"
$ws.Rows.Item(37).RowHeight = 36.5

# --- Row 38: 20 Feb 2020 (Th), 1615-1640 ---
$ws.Range("A38").Value = "20 Feb 2020 (Th)"
$ws.Range("B38").Value = "1615-1640"
$ws.Range("C38").Value = "Harry, Deon, Thuc"
$ws.Range("D38").Value = "Look at homework one last time."
$ws.Range("E38").Value = "Submitted."
$ws.Range("F38").Value = "We spent more time deciding whether the document was good enough than actually working on it."

# Reflection/mood cell with mixed (rich text) formatting:
#   "Need food and water. Badly." -> italic + strikethrough
#   " Went to lunch while waiting for Deon and Thuc to arrive." -> italic only
$ws.Range("G38").Value = "Need food and water. Badly. Went to lunch while waiting for Deon and Thuc to arrive."
$g38 = $ws.Range("G38")
$g38.Characters(1, 27).Font.Italic = $true
$g38.Characters(1, 27).Font.Strikethrough = $true
$g38.Characters(28, 57).Font.Italic = $true
$g38.Font.Italic = $true
$g38.Font.Strikethrough = $true

$ws.Rows.Item(38).RowHeight = 36.5

# --- View state: scrolled/zoomed/selected around the new entries ---
$excel.ActiveWindow.Zoom = 63
$ws.Range("F38").Select() | Out-Null
